$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in row 3 (B3, C3, D3)
$ws.Range("B3").Value = 80
$ws.Range("C3").Value = 75
$ws.Range("D3").Value = 90

# Update selection to B3:D3 with active cell B3
$ws.Range("B3:D3").Select()
